# Add a new "InvalidLogin" worksheet after the existing "ValidLogin" sheet,
# mirroring the username/password header row but with invalid credentials.

$wb = $excel.ActiveWorkbook
$validLogin = $wb.Worksheets.Item("ValidLogin")

# Insert the new sheet right after ValidLogin (Excel defaults to inserting
# before the active sheet, so pass it explicitly as "After").
$invalidLogin = $wb.Worksheets.Add($null, $validLogin)
$invalidLogin.Name = "InvalidLogin"

$invalidLogin.Range("A1").Value = "username"
$invalidLogin.Range("B1").Value = "password"
$invalidLogin.Range("A2").Value = "abcd"
$invalidLogin.Range("B2").Value = "xyz"

# Match the author's final UI state: new sheet active, B2 selected.
$invalidLogin.Activate()
$invalidLogin.Range("B2").Select()
